# Loan RBI, Variable Instalments
#
# The "Repayment schedule" sheet gets a new (blank) column inserted right
# before the existing "Late" column (column N / 14th column), shifting
# "Late", the duplicated "Paid Date" heading and "Outstanding" one column
# to the right (N->O, O->P, P->Q). The newly inserted column keeps the
# same width as the column immediately to its left ("In Advance"), matching
# what Excel does on a manual Insert Column.
#
# The "Repayment schedule" sheet also becomes the active tab/selected sheet
# (it previously was "NewLoanInput"), with a fresh selection at R10.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14th column == "Late").
$ws.Columns.Item(14).Insert()

# Match the width Excel gives a freshly inserted column when it copies the
# formatting of the column to its left ("In Advance", column M).
$ws.Columns.Item(14).ColumnWidth = $ws.Columns.Item(13).ColumnWidth

# Make "Repayment schedule" the active sheet and set the new selection.
$ws.Activate()
$ws.Range("R10").Select()
